# Generate Report for Handoff
# Replaces the old localization-run identifiers/hashes with a new run's,
# bumps a couple of timestamps, clears the (not-yet-produced) handback
# file/hyperlink on the zh-cn and de-de sheets, and flips "Has metadata"
# to True now that the new handoff package carries metadata.

$wb = $excel.ActiveWorkbook

$oldGuid = "31919597-49f0-4bd7-94a7-6977a4835286"
$newGuid = "753f231d-5255-40c6-8a90-9743430aa792"
$newHash = "f8ee534449a79c5911aa99b00394da67f566e8d5"

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "'$newGuid.md"
$wsOverview.Range("G2").Value = "'2017-02-17 09:24:01"

# B2 carries a hyperlink whose cached display text must move to the new
# file name too; drop the old link(s) on the sheet and re-add a clean one
# (Hyperlinks.Add also stamps the cell's text, keeping both in sync).
$wsOverview.Range("A1").Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("B2"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/00c0ca81780699bacd0fbf4d09f6d52e5906f83a/e2e/$newGuid.md",
    [Type]::Missing,
    [Type]::Missing,
    "e2e\$newGuid.md"
)

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("G2").Value = "'$newGuid.$newHash.zh-cn.xlf"
$wsZh.Range("H2").Value = "'2017-02-17 09:23:45"
$wsZh.Range("J2").Value = ""
$wsZh.Range("K2").Value = ""
$wsZh.Range("L2").Value = "'0001-01-01 00:00:00"
$wsZh.Range("Q2").Value = "'True"

$wsZh.Columns.Item(10).ColumnWidth = 18.6506061553955
$wsZh.Columns.Item(11).ColumnWidth = 21.7054767608643

# A2 keeps its hyperlink but the cached display text needs the new file
# name; J2's hyperlink is gone entirely (no handback target yet), so we
# drop every link on the sheet and re-add only the A2 one.
$wsZh.Range("A1").Hyperlinks.Delete()
$wsZh.Hyperlinks.Add(
    $wsZh.Range("A2"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/00c0ca81780699bacd0fbf4d09f6d52e5906f83a/e2e/$newGuid.md",
    [Type]::Missing,
    [Type]::Missing,
    "$newGuid.md"
)

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("G2").Value = "'$newGuid.$newHash.de-de.xlf"
$wsDe.Range("H2").Value = "'2017-02-17 09:24:01"
$wsDe.Range("J2").Value = ""
$wsDe.Range("K2").Value = ""
$wsDe.Range("L2").Value = "'0001-01-01 00:00:00"
$wsDe.Range("Q2").Value = "'True"

$wsDe.Columns.Item(10).ColumnWidth = 18.6506061553955
$wsDe.Columns.Item(11).ColumnWidth = 21.7054767608643

$wsDe.Range("A1").Hyperlinks.Delete()
$wsDe.Hyperlinks.Add(
    $wsDe.Range("A2"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/00c0ca81780699bacd0fbf4d09f6d52e5906f83a/e2e/$newGuid.md",
    [Type]::Missing,
    [Type]::Missing,
    "$newGuid.md"
)
